$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "CCLRC"
$ws.Range("A3").Value = "PPARC"
$ws.Range("A4").Value = "UK Space Agency"
$ws.Range("A5").Value = "British Academy"
$ws.Range("A6").Value = "Royal Academy Eng."
$ws.Range("A7").Value = "STFC"
$ws.Range("A8").Value = "BBSRC"
$ws.Range("A9").Value = "Wellcome"
$ws.Range("A10").Value = "RCUK"
$ws.Range("A11").Value = "NERC"
$ws.Range("A12").Value = "AHRC"
$ws.Range("A13").Value = "MRC"
$ws.Range("A14").Value = "ESRC"
$ws.Range("A15").Value = "Royal Society"
$ws.Range("A16").Value = "EPSRC"
